$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily price table keeps its newest date in row 2 (just below the header)
# and pushes older rows down. Insert a fresh row there for today's update.
$ws.Rows("2:2").Insert()

# Force the date to be stored as plain text (matching every other date cell in
# column A) instead of letting Excel auto-convert "2025-12-06" into a date
# serial number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-06"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# The inserted row copies formatting (bold font/border) from the header row
# above it; strip that back off so the new row matches the plain styling of
# the rest of the data rows.
$ws.Range("A2:D2").ClearFormats()
